$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that should become empty inline strings
$ws.Range("D2").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("I13").Value = ""

# Update cells with new text values
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"
